# Omar Safwat - Week 2 EDA deck edits
# 1) Footer "datetimeFigureOut" field text on the slide master and every
#    slide layout: 03/07/2021 -> 03/11/2021
# 2) Slide 1 title textbox: replace the "Week 2: EDA" / "2021-03-07"
#    paragraphs with four new paragraphs (Name/Location/Team/Date) and
#    grow the textbox to its new auto-fit height.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($j = 1; $j -le $shapes.Count; $j++) {
        $shp = $shapes.Item($j)
        if ($shp.Name -like "Date Placeholder*" -and $shp.HasTextFrame) {
            $tr = $shp.TextFrame.TextRange
            if ($tr.Text -eq "03/07/2021") {
                $tr.Text = "03/11/2021"
            }
        }
    }
}

# -- Part 1: update the footer date field everywhere it lives --
Update-DatePlaceholder($p.SlideMaster.Shapes)

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder($layouts.Item($i).Shapes)
}

# -- Part 2: rewrite the title-slide textbox --
$slide1 = $p.Slides.Item(1)
$titleShape = $slide1.Shapes.Item(2)
$tr = $titleShape.TextFrame.TextRange

# Paragraph 2: "Week 2: EDA" -> "Name: Omar Safwat" (keeps its sz=4000 run)
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "Name: Omar Safwat"

# New paragraph "Location: Egypt" right after it, inheriting sz=4000
$para2.InsertAfter("`rLocation: Egypt") | Out-Null

# New paragraph "Team: Data Science and Analytics" after that one
$para3 = $tr.Paragraphs(3, 1)
$para3.InsertAfter("`rTeam: Data Science and Analytics") | Out-Null

# The original last paragraph ("2021-03-07") is now paragraph 5;
# turn it into the new Date line and drop its old 28pt/bold formatting
# so it matches the other new lines (40pt, not bold).
$para5 = $tr.Paragraphs(5, 1)
$para5.Text = "Date: 11-March-2021"
$para5.Font.Size = 40
$para5.Font.Bold = 0

# The textbox has spAutoFit; set its height explicitly to match the
# grown, five-paragraph content (870857/2380343 offset is unchanged).
$titleShape.Height = 281.1187401574803
